$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.145036666666667
$ws.Range("H2").Value = 3.43511
$ws.Range("I2").Value = 0.4953865629219574
$ws.Range("J2").Value = 0.4953865629219574
$ws.Range("M2").Value = 14.25737566666667
$ws.Range("N2").Value = 42.772127
$ws.Range("O2").Value = 0.2087950866344732
$ws.Range("P2").Value = 0.2087950866344732
$ws.Range("Q2").Value = 16.32521790877445
$ws.Range("R2").Value = 146.92696117897
$ws.Range("S2").Value = 0.103434280322844
$ws.Range("T2").Value = 0.103434280322844
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.145036666666667
$ws.Range("H3").Value = 3.43511
$ws.Range("I3").Value = 0.4953865629219574
$ws.Range("J3").Value = 0.4953865629219574
$ws.Range("N3").Value = 87.128332
$ws.Range("O3").Value = 0.4253229592313036
$ws.Range("P3").Value = 0.4253229592313036
$ws.Range("Q3").Value = 33.25504494850222
$ws.Range("R3").Value = 299.29540453652
$ws.Range("S3").Value = 0.2106992789053913
$ws.Range("T3").Value = 0.2106992789053913
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.145036666666667
$ws.Range("H4").Value = 3.43511
$ws.Range("I4").Value = 0.4953865629219574
$ws.Range("J4").Value = 0.4953865629219574
$ws.Range("M4").Value = 20.11084633333333
$ws.Range("N4").Value = 60.332539
$ws.Range("O4").Value = 0.2945174484164121
$ws.Range("P4").Value = 0.2945174484164122
$ws.Range("Q4").Value = 23.02765644936555
$ws.Range("R4").Value = 207.24890804429
$ws.Range("S4").Value = 0.1458999864915513
$ws.Range("T4").Value = 0.1458999864915513
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.145036666666667
$ws.Range("H5").Value = 3.43511
$ws.Range("I5").Value = 0.4953865629219574
$ws.Range("J5").Value = 0.4953865629219574
$ws.Range("M5").Value = 4.873057999999999
$ws.Range("N5").Value = 14.619174
$ws.Range("O5").Value = 0.07136450571781097
$ws.Range("P5").Value = 0.07136450571781099
$ws.Range("Q5").Value = 5.579830088793333
$ws.Range("R5").Value = 50.21847079913999
$ws.Range("S5").Value = 0.03535301720217075
$ws.Range("T5").Value = 0.03535301720217076
$ws.Range("G6").Value = 0.4713496666666666
$ws.Range("I6").Value = 0.2039238551060172
$ws.Range("J6").Value = 0.2039238551060172
$ws.Range("M6").Value = 14.25737566666667
$ws.Range("N6").Value = 42.772127
$ws.Range("O6").Value = 0.2087950866344732
$ws.Range("P6").Value = 0.2087950866344732
$ws.Range("Q6").Value = 6.720209268024778
$ws.Range("R6").Value = 60.481883412223
$ws.Range("S6").Value = 0.04257829899369663
$ws.Range("T6").Value = 0.04257829899369663
$ws.Range("G7").Value = 0.4713496666666666
$ws.Range("I7").Value = 0.2039238551060172
$ws.Range("J7").Value = 0.2039238551060172
$ws.Range("N7").Value = 87.128332
$ws.Range("O7").Value = 0.4253229592313036
$ws.Range("P7").Value = 0.4253229592313036
$ws.Range("S7").Value = 0.08673349751154683
$ws.Range("T7").Value = 0.08673349751154684
$ws.Range("G8").Value = 0.4713496666666666
$ws.Range("I8").Value = 0.2039238551060172
$ws.Range("J8").Value = 0.2039238551060172
$ws.Range("M8").Value = 20.11084633333333
$ws.Range("N8").Value = 60.332539
$ws.Range("O8").Value = 0.2945174484164121
$ws.Range("P8").Value = 0.2945174484164122
$ws.Range("Q8").Value = 9.47924071560122
$ws.Range("R8").Value = 85.313166440411
$ws.Range("S8").Value = 0.06005913347706233
$ws.Range("T8").Value = 0.06005913347706235
$ws.Range("G9").Value = 0.4713496666666666
$ws.Range("I9").Value = 0.2039238551060172
$ws.Range("J9").Value = 0.2039238551060172
$ws.Range("M9").Value = 4.873057999999999
$ws.Range("N9").Value = 14.619174
$ws.Range("O9").Value = 0.07136450571781097
$ws.Range("P9").Value = 0.07136450571781099
$ws.Range("Q9").Value = 2.296914263947333
$ws.Range("R9").Value = 20.672228375526
$ws.Range("S9").Value = 0.01455292512371142
$ws.Range("T9").Value = 0.01455292512371143
$ws.Range("G10").Value = 0.6323219999999999
$ws.Range("H10").Value = 1.896966
$ws.Range("I10").Value = 0.2735666300991275
$ws.Range("J10").Value = 0.2735666300991275
$ws.Range("M10").Value = 14.25737566666667
$ws.Range("N10").Value = 42.772127
$ws.Range("O10").Value = 0.2087950866344732
$ws.Range("P10").Value = 0.2087950866344732
$ws.Range("Q10").Value = 9.015252296298
$ws.Range("R10").Value = 81.137270666682
$ws.Range("S10").Value = 0.0571193682318482
$ws.Range("T10").Value = 0.05711936823184822
$ws.Range("G11").Value = 0.6323219999999999
$ws.Range("H11").Value = 1.896966
$ws.Range("I11").Value = 0.2735666300991275
$ws.Range("J11").Value = 0.2735666300991275
$ws.Range("N11").Value = 87.128332
$ws.Range("O11").Value = 0.4253229592313036
$ws.Range("P11").Value = 0.4253229592313036
$ws.Range("Q11").Value = 18.364387048968
$ws.Range("R11").Value = 165.279483440712
$ws.Range("S11").Value = 0.1163541686606963
$ws.Range("T11").Value = 0.1163541686606963
$ws.Range("G12").Value = 0.6323219999999999
$ws.Range("H12").Value = 1.896966
$ws.Range("I12").Value = 0.2735666300991275
$ws.Range("J12").Value = 0.2735666300991275
$ws.Range("M12").Value = 20.11084633333333
$ws.Range("N12").Value = 60.332539
$ws.Range("O12").Value = 0.2945174484164121
$ws.Range("P12").Value = 0.2945174484164122
$ws.Range("Q12").Value = 12.716530575186
$ws.Range("R12").Value = 114.448775176674
$ws.Range("S12").Value = 0.08057014586867148
$ws.Range("T12").Value = 0.08057014586867151
$ws.Range("G13").Value = 0.6323219999999999
$ws.Range("H13").Value = 1.896966
$ws.Range("I13").Value = 0.2735666300991275
$ws.Range("J13").Value = 0.2735666300991275
$ws.Range("M13").Value = 4.873057999999999
$ws.Range("N13").Value = 14.619174
$ws.Range("O13").Value = 0.07136450571781097
$ws.Range("P13").Value = 0.07136450571781099
$ws.Range("Q13").Value = 3.081341780676
$ws.Range("R13").Value = 27.732076026084
$ws.Range("S13").Value = 0.01952294733791146
$ws.Range("T13").Value = 0.01952294733791147
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.062692
$ws.Range("H14").Value = 0.188076
$ws.Range("I14").Value = 0.02712295187289783
$ws.Range("J14").Value = 0.02712295187289783
$ws.Range("M14").Value = 14.25737566666667
$ws.Range("N14").Value = 42.772127
$ws.Range("O14").Value = 0.2087950866344732
$ws.Range("P14").Value = 0.2087950866344732
$ws.Range("Q14").Value = 0.8938233952946667
$ws.Range("R14").Value = 8.044410557652
$ws.Range("S14").Value = 0.005663139086084349
$ws.Range("T14").Value = 0.00566313908608435
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.062692
$ws.Range("H15").Value = 0.188076
$ws.Range("I15").Value = 0.02712295187289783
$ws.Range("J15").Value = 0.02712295187289783
$ws.Range("N15").Value = 87.128332
$ws.Range("O15").Value = 0.4253229592313036
$ws.Range("P15").Value = 0.4253229592313036
$ws.Range("Q15").Value = 1.820749796581333
$ws.Range("R15").Value = 16.386748169232
$ws.Range("S15").Value = 0.01153601415366913
$ws.Range("T15").Value = 0.01153601415366913
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.062692
$ws.Range("H16").Value = 0.188076
$ws.Range("I16").Value = 0.02712295187289783
$ws.Range("J16").Value = 0.02712295187289783
$ws.Range("M16").Value = 20.11084633333333
$ws.Range("N16").Value = 60.332539
$ws.Range("O16").Value = 0.2945174484164121
$ws.Range("P16").Value = 0.2945174484164122
$ws.Range("Q16").Value = 1.260789178329333
$ws.Range("R16").Value = 11.347102604964
$ws.Range("S16").Value = 0.007988182579127015
$ws.Range("T16").Value = 0.007988182579127017
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.062692
$ws.Range("H17").Value = 0.188076
$ws.Range("I17").Value = 0.02712295187289783
$ws.Range("J17").Value = 0.02712295187289783
$ws.Range("M17").Value = 4.873057999999999
$ws.Range("N17").Value = 14.619174
$ws.Range("O17").Value = 0.07136450571781097
$ws.Range("P17").Value = 0.07136450571781099
$ws.Range("Q17").Value = 0.3055017521359999
$ws.Range("R17").Value = 2.749515769224
$ws.Range("S17").Value = 0.001935616054017329
$ws.Range("T17").Value = 0.001935616054017329
